$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (DAMSLTag, DialogAct) updates, derived from the diff
$updates = @(
    @{ Row = 5; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 18; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 23; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 28; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 36; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 38; I = '%'; J = 'Uninterpretable' }
    @{ Row = 44; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 46; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 54; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 65; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 73; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 74; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 84; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 85; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 97; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 106; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 121; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 126; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 131; I = 'ba'; J = 'Appreciation' }
    @{ Row = 132; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 139; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 150; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 152; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 154; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 169; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 170; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 171; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 178; I = '%'; J = 'Uninterpretable' }
    @{ Row = 181; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 183; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 199; I = 'ba'; J = 'Appreciation' }
    @{ Row = 207; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 212; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 213; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 220; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 228; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 231; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 247; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 250; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 253; I = '%'; J = 'Uninterpretable' }
    @{ Row = 254; I = '%'; J = 'Uninterpretable' }
    @{ Row = 259; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 269; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 272; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 277; I = '%'; J = 'Uninterpretable' }
    @{ Row = 292; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 312; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 335; I = 'ba'; J = 'Appreciation' }
    @{ Row = 355; I = 'ba'; J = 'Appreciation' }
    @{ Row = 359; I = 'ba'; J = 'Appreciation' }
    @{ Row = 361; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 363; I = 'ba'; J = 'Appreciation' }
    @{ Row = 370; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 384; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 395; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 397; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 402; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 406; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 408; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 410; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 411; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 422; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 423; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 427; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 431; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 433; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 447; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 448; I = 'ba'; J = 'Appreciation' }
    @{ Row = 449; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 451; I = '%'; J = 'Uninterpretable' }
    @{ Row = 453; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 457; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 461; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 462; I = 'ba'; J = 'Appreciation' }
    @{ Row = 469; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 470; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 483; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 486; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 487; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 490; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 495; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 499; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 502; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 511; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 523; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 526; I = 'ba'; J = 'Appreciation' }
    @{ Row = 528; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 533; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 539; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 542; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 544; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 552; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 553; I = 'ba'; J = 'Appreciation' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
